# Update the CauHoi sheet: every question's "maKhoi" (grade) value moves
# from K10 to K11 (sua file cau hoi Ly11.xlsx).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CauHoi")

$ws.Range("E2:E25").Value = "K11"

# Match the resulting selection left behind by the edit.
$ws.Range("E3:E25").Select() | Out-Null
